$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.785.10'
$ws.Range('E2').Value = '  +2.69%  '
$ws.Range('D3').Value = '3.796.62'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '701.26'
$ws.Range('E5').Value = '  +9.02%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '172.64'
$ws.Range('E6').Value = '  +4.40%  '
$ws.Range('D7').Value = '3.795.18'
$ws.Range('E7').Value = '  +0.85%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.527'
$ws.Range('E9').Value = '  +0.94%  '
$ws.Range('E10').Value = '  +2.62%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.38'
$ws.Range('E11').Value = '  +6.60%  '
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('E13').Value = '  +8.30%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '36.43'
$ws.Range('E14').Value = '  +4.47%  '
$ws.Range('D15').Value = '4.436.42'
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('D16').Value = '3.835.47'
$ws.Range('E16').Value = '  +1.70%  '
$ws.Range('D17').Value = '70.799.73'
$ws.Range('E17').Value = '  +2.73%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '17.87'
$ws.Range('E18').Value = '  +1.28%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.20'
$ws.Range('E19').Value = '  +2.87%  '
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.08'
$ws.Range('E21').Value = '  +16.03%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '482.02'
$ws.Range('E22').Value = '  +2.04%  '
$ws.Range('E23').Value = '  +1.43%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '84.26'
$ws.Range('E24').Value = '  +3.16%  '
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.44'
$ws.Range('E26').Value = '  +2.42%  '
$ws.Range('E27').Value = '  +3.70%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.48'
$ws.Range('E28').Value = '  +4.24%  '
$ws.Range('D29').Value = '3.947.00'
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.12'
$ws.Range('E31').Value = '  +16.05%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.53'
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '29.50'
$ws.Range('E34').Value = '  +3.43%  '
$ws.Range('E35').Value = '  +4.70%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '9.23'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  +2.28%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.44'
$ws.Range('E39').Value = '  +6.14%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.08'
$ws.Range('E40').Value = '  +5.37%  '
$ws.Range('E41').Value = '  +12.09%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.972'
$ws.Range('E42').Value = '  +1.51%  '
$ws.Range('E43').Value = '  +22.59%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('B46').Value = 'Arweave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '45.80'
$ws.Range('E46').Value = '  +2.09%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '162.30'
$ws.Range('E47').Value = '  +4.24%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '48.75'
$ws.Range('E48').Value = '  +2.30%  '
$ws.Range('E49').Value = '  +3.02%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.57'
$ws.Range('E51').Value = '  +2.57%  '
